# Update "Lista de Riscos" — revise probability estimates for a few risks
# and move the live selection to where the user left off reviewing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos")

# Risk #2 ("Desistência de membro da equipe"): probability re-assessed to 0%
$ws.Range("G4").Value = 0

# Risk #3 ("Sobreposição de atividades"): probability re-assessed to 25%
$ws.Range("G5").Value = 0.25

# Risk #7 ("Alteração de plataforma do projeto"): probability re-assessed to 0%
$ws.Range("G9").Value = 0

# Leave the selection where the user was last working
$ws.Range("G14").Select()
